$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header column (H1), reusing the existing header style
# (copy formats from the neighboring header cell so no new style is minted)
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").Value = "Save"

# Add the corresponding data cell (H2)
$ws.Range("H2").Value = 0

$excel.CutCopyMode = $false
